$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 38, shifting the existing rows 38-53 (飛龍 sub-entries) down to 39-54.
# The new row inherits formatting/style from the row above (row 37).
$ws.Rows("38:38").Insert()

# Fill A38 (new row) and C38 (new row) with the new monster entry.
# Also fill in column A for rows 39-54 (previously blank, now carrying the "飛龍" category)
# and append new rows 55-75 for the remaining monsters.
$ws.Range("A38").Value = "飛龍"
$ws.Range("C38").Value = "霜刃冰牙龍"
$ws.Range("A39").Value = "飛龍"
$ws.Range("C39").Value = "爆麟龍"
$ws.Range("A40").Value = "飛龍"
$ws.Range("C40").Value = "紅蓮爆麟龍"
$ws.Range("A41").Value = "飛龍"
$ws.Range("C41").Value = "角龍"
$ws.Range("A42").Value = "飛龍"
$ws.Range("C42").Value = "黑角龍"
$ws.Range("A43").Value = "飛龍"
$ws.Range("C43").Value = "風漂龍"
$ws.Range("A44").Value = "飛龍"
$ws.Range("C44").Value = "霜翼風漂龍"
$ws.Range("A45").Value = "飛龍"
$ws.Range("C45").Value = "迅龍"
$ws.Range("A46").Value = "飛龍"
$ws.Range("C46").Value = "浮空龍"
$ws.Range("A47").Value = "飛龍"
$ws.Range("C47").Value = "浮眠龍"
$ws.Range("A48").Value = "飛龍"
$ws.Range("C48").Value = "火龍"
$ws.Range("A49").Value = "飛龍"
$ws.Range("C49").Value = "蒼火龍"
$ws.Range("A50").Value = "飛龍"
$ws.Range("C50").Value = "雌火龍"
$ws.Range("A51").Value = "飛龍"
$ws.Range("C51").Value = "櫻火龍"
$ws.Range("A52").Value = "飛龍"
$ws.Range("C52").Value = "轟龍"
$ws.Range("A53").Value = "飛龍"
$ws.Range("C53").Value = "黑轟龍"
$ws.Range("A54").Value = "飛龍"
$ws.Range("C54").Value = "銀火龍"
$ws.Range("A55").Value = "飛龍"
$ws.Range("C55").Value = "金火龍"
$ws.Range("A56").Value = "魚龍"
$ws.Range("C56").Value = "冰魚龍"
$ws.Range("A57").Value = "魚龍"
$ws.Range("C57").Value = "泥魚龍"
$ws.Range("A58").Value = "魚龍"
$ws.Range("C58").Value = "熔岩龍"
$ws.Range("A59").Value = "古龍"
$ws.Range("C59").Value = "貝西摩斯"
$ws.Range("A60").Value = "古龍"
$ws.Range("C60").Value = "剛龍"
$ws.Range("A61").Value = "古龍"
$ws.Range("C61").Value = "冰咒龍"
$ws.Range("A62").Value = "古龍"
$ws.Range("C62").Value = "滅盡龍"
$ws.Range("A63").Value = "古龍"
$ws.Range("C63").Value = "殲世滅盡龍"
$ws.Range("A64").Value = "古龍"
$ws.Range("C64").Value = "屍套龍"
$ws.Range("A65").Value = "古龍"
$ws.Range("C65").Value = "霧瘴屍套龍"
$ws.Range("A66").Value = "古龍"
$ws.Range("C66").Value = "炎王龍"
$ws.Range("A67").Value = "古龍"
$ws.Range("C67").Value = "炎妃龍"
$ws.Range("A68").Value = "古龍"
$ws.Range("C68").Value = "麒麟"
$ws.Range("A69").Value = "古龍"
$ws.Range("C69").Value = "冥燈龍"
$ws.Range("A70").Value = "古龍"
$ws.Range("C70").Value = "溟波龍"
$ws.Range("A71").Value = "古龍"
$ws.Range("C71").Value = "天地煌啼龍"
$ws.Range("A72").Value = "古龍"
$ws.Range("C72").Value = "絢輝龍"
$ws.Range("A73").Value = "古龍"
$ws.Range("C73").Value = "冥赤龍"
$ws.Range("A74").Value = "古龍"
$ws.Range("C74").Value = "煌黑龍"
$ws.Range("A75").Value = "古龍"
$ws.Range("C75").Value = "黑龍"

# Update the view selection to match the new active cell.
$ws.Range("D63").Select()
